# Carbonell_Cathleen_ProblemSolving.docx -- "Predicting Fingers" problem:
# add the worked "Solution" write-up for Part A (finding which finger the
# count of 10 lands on), including the numbered list of fingers.
#
# Before:
#   ... Define. The goal is to predict ... 1000.
#   <empty paragraph, carries the _GoBack bookmark>
#   Solution. One suggestion ... tedious.<br/>
#   But since there is a pattern to her method-- and one can draft it out as such:
#   <empty NormalWeb paragraph>
#
# After:
#   ... Define. The goal is to predict ... 1000.
#   <empty paragraph (bookmark removed from here)>
#   Solution. One suggestion ... tedious.<br/>
#   But since there is a pattern to her method-- and one can draft it out as such:
#   <empty paragraph>
#   Final Resolve part a: Finding 10<br/><br/>
#   There are 5 digits on one hand. The first set of 5, numbers to digits, is as follows:
#   - Thumb / Index finger / Middle finger / Ring finger / Pinky
#   <empty ListParagraph-styled paragraph>
#   Then, she goes back and counts her ring finger as 6. If this follows, to 10, the pattern is as follows:
#   - Ring / Middle / Index / Thumb / " Index"
#   <empty paragraph>
#   <NormalWeb paragraph, now carrying the _GoBack bookmark>

$d = $word.ActiveDocument

function Wrap-WordXmlBody([string]$bodyXml) {
    return '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' +
        $bodyXml +
        '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
}

# --- Step 1 -----------------------------------------------------------
# The paragraph right after "Define. The goal is to predict..." is an
# otherwise-empty paragraph that happens to carry the (hidden) _GoBack
# bookmark. Locate it via the bookmark itself, rather than a hard-coded
# paragraph index, and clear it down to a plain empty paragraph -- the
# bookmark is re-homed onto the very last paragraph in step 2.
$goBack = $d.Bookmarks.Item("_GoBack")
$bookmarkPara = $goBack.Range.Paragraphs.Item(1)
$bookmarkPara.Range.InsertXML((Wrap-WordXmlBody "<w:p/>"))

# --- Step 2 -----------------------------------------------------------
# Replace everything from the "Solution." paragraph through to the final
# (empty, NormalWeb-styled) paragraph of the document with the kept
# "Solution"/"But since..." paragraphs plus the brand-new write-up.
$solutionRange = $d.Content
$solutionRange.Find.Execute("Solution.", $true) | Out-Null
$solutionPara = $solutionRange.Paragraphs.Item(1)

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)

$replaceRange = $d.Range($solutionPara.Range.Start, $lastPara.Range.End)

$newBodyXml = @'
<w:p><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">Solution. </w:t></w:r><w:r><w:t xml:space="preserve">One suggestion to solve this is that we </w:t></w:r><w:r><w:rPr><w:i/></w:rPr><w:t>could</w:t></w:r><w:r><w:t xml:space="preserve"> count to 10, 100, and 1000 ourselves the way the girl does, which would be extremely tedious.</w:t></w:r><w:r><w:br/></w:r></w:p><w:p><w:r><w:t>But since there is a pattern to her method</w:t></w:r></w:p><w:p/><w:p><w:r><w:rPr><w:b/></w:rPr><w:t>Final Resolve part a: Finding 10</w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:br/></w:r><w:r><w:rPr><w:b/></w:rPr><w:br/></w:r><w:r><w:t>There are 5 digits on one hand. The first set of 5, numbers to digits, is as follows:</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="13"/></w:numPr></w:pPr><w:r><w:t>Thumb</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="13"/></w:numPr></w:pPr><w:r><w:t>Index finger</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="13"/></w:numPr></w:pPr><w:r><w:t>Middle finger</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="13"/></w:numPr></w:pPr><w:r><w:t>Ring finger</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="13"/></w:numPr></w:pPr><w:r><w:t>Pinky</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:ind w:left="1080"/></w:pPr></w:p><w:p><w:r><w:t>Then, she goes back and counts her ring finger as 6. If this follows, to 10, the pattern is as follows:</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="13"/></w:numPr></w:pPr><w:r><w:t>Ring</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="13"/></w:numPr></w:pPr><w:r><w:t>Middle</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="13"/></w:numPr></w:pPr><w:r><w:t>Index</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="13"/></w:numPr></w:pPr><w:r><w:t>Thumb</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="13"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve"> Index</w:t></w:r></w:p><w:p/><w:p><w:pPr><w:pStyle w:val="NormalWeb"/></w:pPr><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>
'@

# The em dash in "her method<EM DASH> and one can draft..." doesn't type
# cleanly through the here-string above on every host code page, so splice
# it in as a literal unicode character from PowerShell itself.
$emDash = [char]0x2014
$newBodyXml = $newBodyXml.Replace("her method</w:t>", "her method" + $emDash + " and one can draft it out as such:</w:t>")

$replaceRange.InsertXML((Wrap-WordXmlBody $newBodyXml))
